$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New progress entries: date (serial), completed chapters.
# Percentage column keeps the existing "B/200*100" pattern from rows 2-7.
$newRows = @(
    @{ Row = 8;  Date = 45789; Chapters = 51 },
    @{ Row = 9;  Date = 45790; Chapters = 53 },
    @{ Row = 10; Date = 45793; Chapters = 58 }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Date
    $ws.Cells.Item($r, 1).NumberFormat = "d-mmm"
    $ws.Cells.Item($r, 2).Value = $entry.Chapters
    $ws.Cells.Item($r, 3).Formula = "=B" + $r + "/200*100"
}

$ws.Range("A11").Select()
